$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the "mtitleStyle" named style (same as used on row 9 header labels)
# to the row-label cells A10, A11, A12.
$ws.Range("A10").Style = "mtitleStyle"
$ws.Range("A11").Style = "mtitleStyle"
$ws.Range("A12").Style = "mtitleStyle"

# Update the "Wrong" marking value from -3 to -1
$ws.Range("C11").Value = "-1"

# Update the computed "Wrong" total from -6 to -2
$ws.Range("C12").Value = -2

# Update the total score text to reflect the corrected marking
$ws.Range("E12").Value = "128/140"
